# coverpage_template.docx fixes:
#  1. Remove the left indent on the "Rajshahi University ..." title line
#     (keep the existing right indent).
#  2. Shrink the four "Course code / Course Name / Report Name /
#     Date of Submission" lines from 16pt to 14pt (both sz and szCs / the
#     complex-script size).
#  3. Drop the stray empty (18pt) paragraph that used to sit between the
#     "Date of Submission" line and the following blank (14pt) paragraph
#     right before the submitted-by/submitted-to table.
#  4. Re-indent that table and resize its two columns.
#  5. Restore normal 1-inch page margins (they were all zeroed out).
#  6. Strip the picture watermark out of all three headers.

$d = $word.ActiveDocument

# 1. "Rajshahi University of Engineering & Technology" title paragraph.
$titlePara = $d.Paragraphs.Item(4)
$titlePara.Format.LeftIndent = 0

# 2. Course code / Course Name / Report Name / Date of Submission lines.
for ($i = 15; $i -le 18; $i++) {
    $rng = $d.Paragraphs.Item($i).Range
    $rng.Font.Size = 14
    $rng.Font.SizeBi = 14
}

# 3. Delete the leftover empty 18pt paragraph (directly after the
#    "Date of Submission" line, before the blank 14pt paragraph).
$d.Paragraphs.Item(19).Range.Delete()

# 4. Table indentation + column widths.
$tbl = $d.Tables.Item(1)
$tbl.Rows.LeftIndent = 8.75
$tbl.Cell(1, 1).Width = 229.5
$tbl.Cell(1, 2).Width = 225

# 5. Page margins back to 1 inch all round.
$ps = $d.Sections.Item(1).PageSetup
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.LeftMargin = 72
$ps.RightMargin = 72

# 6. Remove the watermark picture from each of the 3 headers
#    (even, first, default/primary all carry one copy of it).
$sec = $d.Sections.Item(1)
$headerTypes = @(1, 2, 3)
foreach ($ht in $headerTypes) {
    $hdr = $sec.Headers.Item($ht)
    while ($hdr.Shapes.Count -gt 0) {
        $hdr.Shapes.Item(1).Delete()
    }
}
